$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Förändrad" (last-changed) date in column C for all data rows (2-176): 45183 -> 45184 ---
for ($r = 2; $r -le 176; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45184
}

# --- Rebuild hyperlink formulas for rows 2-5 (species-hit rows) ---
# Row 2: A 30840-2023
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/artfynd/A 30840-2023.xlsx, "A 30840-2023"")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/kartor/A 30840-2023.png", "A 30840-2023")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/knärot/A 30840-2023.png", "A 30840-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/klagomål/A 30840-2023.docx", "A 30840-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/klagomålsmail/A 30840-2023.docx", "A 30840-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/tillsyn/A 30840-2023.docx", "A 30840-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/tillsynsmail/A 30840-2023.docx", "A 30840-2023")'

# Row 3: A 30841-2023
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/artfynd/A 30841-2023.xlsx, "A 30841-2023"")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/kartor/A 30841-2023.png", "A 30841-2023")'
$ws.Range("U3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/knärot/A 30841-2023.png", "A 30841-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/klagomål/A 30841-2023.docx", "A 30841-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/klagomålsmail/A 30841-2023.docx", "A 30841-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/tillsyn/A 30841-2023.docx", "A 30841-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/tillsynsmail/A 30841-2023.docx", "A 30841-2023")'

# Row 4: A 30839-2023
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/artfynd/A 30839-2023.xlsx, "A 30839-2023"")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/kartor/A 30839-2023.png", "A 30839-2023")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/klagomål/A 30839-2023.docx", "A 30839-2023")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/klagomålsmail/A 30839-2023.docx", "A 30839-2023")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/tillsyn/A 30839-2023.docx", "A 30839-2023")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/tillsynsmail/A 30839-2023.docx", "A 30839-2023")'

# Row 5: A 33036-2023
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SOLLEFTEA/artfynd/A 33036-2023.xlsx, "A 33036-2023"")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SOLLEFTEA/kartor/A 33036-2023.png", "A 33036-2023")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SOLLEFTEA/klagomål/A 33036-2023.docx", "A 33036-2023")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SOLLEFTEA/klagomålsmail/A 33036-2023.docx", "A 33036-2023")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SOLLEFTEA/tillsyn/A 33036-2023.docx", "A 33036-2023")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_SOLLEFTEA/tillsynsmail/A 33036-2023.docx", "A 33036-2023")'
